$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 213
$ws.Cells.Item(213, 1).NumberFormat = "@"
$ws.Cells.Item(213, 1).Value = '2025-10-03'
$ws.Cells.Item(213, 2).Value = 'Süper Lig'
$ws.Cells.Item(213, 3).Value = 'Trabzonspor'
$ws.Cells.Item(213, 4).Value = 'Kayserispor'
$ws.Cells.Item(213, 5).Value = 'Home Win'
$ws.Cells.Item(213, 6).NumberFormat = "@"
$ws.Cells.Item(213, 6).Value = '85.01%'
$ws.Cells.Item(213, 7).Value = 1.5
$ws.Cells.Item(213, 8).NumberFormat = "@"
$ws.Cells.Item(213, 8).Value = '26.24%'
$ws.Cells.Item(213, 9).Value = 3.2
$ws.Cells.Item(213, 10).Value = 0.05
$ws.Cells.Item(213, 11).Value = 0.5502407679019199
$ws.Cells.Item(213, 12).Value = 'Pending'

# Row 214
$ws.Cells.Item(214, 1).NumberFormat = "@"
$ws.Cells.Item(214, 1).Value = '2025-10-03'
$ws.Cells.Item(214, 2).Value = 'Bundesliga'
$ws.Cells.Item(214, 3).Value = '1899 Hoffenheim'
$ws.Cells.Item(214, 4).Value = '1.FC Köln'
$ws.Cells.Item(214, 5).Value = 'Home Win'
$ws.Cells.Item(214, 6).NumberFormat = "@"
$ws.Cells.Item(214, 6).Value = '65.09%'
$ws.Cells.Item(214, 7).Value = 1.9
$ws.Cells.Item(214, 8).NumberFormat = "@"
$ws.Cells.Item(214, 8).Value = '22.43%'
$ws.Cells.Item(214, 9).Value = 1.7
$ws.Cells.Item(214, 10).Value = 0.02629742390973933
$ws.Cells.Item(214, 11).Value = 0.2629742390973933
$ws.Cells.Item(214, 12).Value = 'Pending'

# Row 215
$ws.Cells.Item(215, 1).NumberFormat = "@"
$ws.Cells.Item(215, 1).Value = '2025-10-03'
$ws.Cells.Item(215, 2).Value = 'Jupiler Pro League'
$ws.Cells.Item(215, 3).Value = 'Gent'
$ws.Cells.Item(215, 4).Value = 'Charleroi'
$ws.Cells.Item(215, 5).Value = 'Home Win'
$ws.Cells.Item(215, 6).NumberFormat = "@"
$ws.Cells.Item(215, 6).Value = '57.63%'
$ws.Cells.Item(215, 7).Value = 2.05
$ws.Cells.Item(215, 8).NumberFormat = "@"
$ws.Cells.Item(215, 8).Value = '16.96%'
$ws.Cells.Item(215, 9).Value = 1.1
$ws.Cells.Item(215, 10).Value = 0.01727445339087996
$ws.Cells.Item(215, 11).Value = 0.1727445339087995
$ws.Cells.Item(215, 12).Value = 'Pending'

# Row 216
$ws.Cells.Item(216, 1).NumberFormat = "@"
$ws.Cells.Item(216, 1).Value = '2025-10-03'
$ws.Cells.Item(216, 2).Value = 'Ligue 1'
$ws.Cells.Item(216, 3).Value = 'Paris FC'
$ws.Cells.Item(216, 4).Value = 'Lorient'
$ws.Cells.Item(216, 5).Value = 'Home Win'
$ws.Cells.Item(216, 6).NumberFormat = "@"
$ws.Cells.Item(216, 6).Value = '70.16%'
$ws.Cells.Item(216, 7).Value = 1.83
$ws.Cells.Item(216, 8).NumberFormat = "@"
$ws.Cells.Item(216, 8).Value = '27.11%'
$ws.Cells.Item(216, 9).Value = 2.2
$ws.Cells.Item(216, 10).Value = 0.03420866828729644
$ws.Cells.Item(216, 11).Value = 0.3420866828729644
$ws.Cells.Item(216, 12).Value = 'Pending'

# Row 217
$ws.Cells.Item(217, 1).NumberFormat = "@"
$ws.Cells.Item(217, 1).Value = '2025-10-03'
$ws.Cells.Item(217, 2).Value = 'Premier League'
$ws.Cells.Item(217, 3).Value = 'Bournemouth'
$ws.Cells.Item(217, 4).Value = 'Fulham'
$ws.Cells.Item(217, 5).Value = 'Home Win'
$ws.Cells.Item(217, 6).NumberFormat = "@"
$ws.Cells.Item(217, 6).Value = '68.87%'
$ws.Cells.Item(217, 7).Value = 1.85
$ws.Cells.Item(217, 8).NumberFormat = "@"
$ws.Cells.Item(217, 8).Value = '26.13%'
$ws.Cells.Item(217, 9).Value = 2
$ws.Cells.Item(217, 10).Value = 0.03223787987903293
$ws.Cells.Item(217, 11).Value = 0.3223787987903293
$ws.Cells.Item(217, 12).Value = 'Pending'

# Row 218
$ws.Cells.Item(218, 1).NumberFormat = "@"
$ws.Cells.Item(218, 1).Value = '2025-10-03'
$ws.Cells.Item(218, 2).Value = 'La Liga'
$ws.Cells.Item(218, 3).Value = 'Osasuna'
$ws.Cells.Item(218, 4).Value = 'Getafe'
$ws.Cells.Item(218, 5).Value = 'Home Win'
$ws.Cells.Item(218, 6).NumberFormat = "@"
$ws.Cells.Item(218, 6).Value = '52.29%'
$ws.Cells.Item(218, 7).Value = 2.2
$ws.Cells.Item(218, 8).NumberFormat = "@"
$ws.Cells.Item(218, 8).Value = '13.90%'
$ws.Cells.Item(218, 9).Value = 0.8
$ws.Cells.Item(218, 10).Value = 0.01253911440951226
$ws.Cells.Item(218, 11).Value = 0.1253911440951226
$ws.Cells.Item(218, 12).Value = 'Pending'

# Row 219
$ws.Cells.Item(219, 1).NumberFormat = "@"
$ws.Cells.Item(219, 1).Value = '2025-10-03'
$ws.Cells.Item(219, 2).Value = 'Liga de Expansión MX'
$ws.Cells.Item(219, 3).Value = 'Tapatío'
$ws.Cells.Item(219, 4).Value = 'Alebrijes de Oaxaca'
$ws.Cells.Item(219, 5).Value = 'Home Win'
$ws.Cells.Item(219, 6).NumberFormat = "@"
$ws.Cells.Item(219, 6).Value = '79.54%'
$ws.Cells.Item(219, 7).Value = 1.62
$ws.Cells.Item(219, 8).NumberFormat = "@"
$ws.Cells.Item(219, 8).Value = '27.57%'
$ws.Cells.Item(219, 9).Value = 2.9
$ws.Cells.Item(219, 10).Value = 0.04654222003881247
$ws.Cells.Item(219, 11).Value = 0.4654222003881247
$ws.Cells.Item(219, 12).Value = 'Pending'
